$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "329.01"
    "D3" = "40.21"
    "E3" = "6.95%"
    "D4" = "5.272"
    "E4" = "2.17%"
    "D5" = "0.08099"
    "E5" = "2.18%"
    "D6" = "4.528"
    "E6" = "2.22%"
    "E7" = "4.24%"
    "D8" = "1.917"
    "E8" = "-0.62%"
    "E9" = "-1.40%"
    "D10" = "0.9350"
    "E10" = "0.37%"
    "D11" = "0.1354"
    "E11" = "23.73%"
    "D12" = "0.1963"
    "E12" = "2.12%"
    "D13" = "0.09325"
    "E13" = "2.82%"
    "D14" = "0.03438"
    "E14" = "4.04%"
    "D15" = "0.09551"
    "E15" = "-0.47%"
    "D16" = "0.001392"
    "E16" = "-0.04%"
    "D17" = "0.006086"
    "E17" = "5.92%"
    "D18" = "3.362"
    "E18" = "-6.66%"
    "D20" = "7.288"
    "E20" = "21.49%"
    "E21" = "1.96%"
    "D22" = "0.2310"
    "E22" = "-10.85%"
    "D23" = "0.04439"
    "E23" = "0.62%"
    "E24" = "-1.08%"
    "E25" = "-5.81%"
    "E26" = "-5.21%"
    "D27" = "0.0003988"
    "E27" = "-0.13%"
    "D39" = "0.02512"
    "E39" = "11.81%"
    "D40" = "0.05241"
    "E40" = "2.65%"
    "D41" = "0.007621"
    "E41" = "2.26%"
    "E42" = "5.70%"
    "D43" = "0.008663"
    "E43" = "-3.77%"
    "E44" = "1.35%"
    "D45" = "0.008135"
    "D46" = "0.00006655"
    "E46" = "0.06%"
    "E47" = "-0.10%"
    "E48" = "148.06%"
    "D49" = "0.002893"
    "E49" = "1.13%"
    "D50" = "0.00002099"
    "E50" = "-0.10%"
    "D51" = "0.0001999"
    "E51" = "-0.10%"
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.ClearFormats()
}
